$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a new "Also, all filenames are driven by the YAML configuraiton
#    file." paragraph (BodyText style) right after the existing
#    "All writing is driven by the YAML configuration file..." paragraph.
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*All writing is driven by the YAML configuration file*") {
        # Create the structural paragraph break (inherits the BodyText style
        # from the paragraph it splits off from).
        $p.Range.InsertParagraphAfter()
        $newPara = $p.Next()

        # Give the brand-new (empty) paragraph a placeholder character so its
        # Range is non-empty; this lets the following InsertXML perform a
        # genuine "replace range contents" instead of a zero-width insert,
        # which keeps the output free of stray empty runs/paragraphs.
        $newPara.Range.Text = "X"
        $newPara2 = $p.Next()
        $target = $newPara2.Range
        $targetNoMark = $d.Range($target.Start, $target.End - 1)

        $flatOpc = '<?xml version="1.0" standalone="yes"?>' +
            '<?mso-application progid="Word.Document"?>' +
            '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
            '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
            '<pkg:xmlData>' +
            '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
            '<w:body><w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t xml:space="preserve">Also, all filenames are driven by the YAML configuraiton file.</w:t></w:r></w:p></w:body>' +
            '</w:document>' +
            '</pkg:xmlData></pkg:part></pkg:package>'

        $targetNoMark.InsertXML($flatOpc)
        break
    }
}

# ---------------------------------------------------------------------------
# 2) Update the sample Windows paths used in the Syntax section so they point
#    at C:\dev\... instead of C:\Users\user\Desktop\....
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "        pyNewCode C:\Users\user\Desktop\MyProject", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "        pyNewCode C:\dev\MyProject", 2)

$d.Content.Find.Execute(
    "        (will build the new in C:\Users\user\Desktop\myProject directory)", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "        (will build the new project in C:\dev\myProject directory)", 2)

$d.Content.Find.Execute(
    "        python pyNewCode.py C:\Users\user\Desktop\MyProject", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "        python pyNewCode.py C:\dev\MyProject", 2)

$d.Content.Find.Execute(
    "        (will build the new project in C:\Users\user\Desktop\MyProject directory)", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "        (will build the new project in C:\dev\MyProject directory)", 2)

# ---------------------------------------------------------------------------
# 3) Bold the standalone "todo.md file" and "project.yaml file" captions
#    (the FirstParagraph-styled lines introducing each example's source
#    block), using InsertXML so the run gets both <w:b/> and <w:bCs/>
#    exactly as produced by Word, without touching the paragraph mark.
# ---------------------------------------------------------------------------
function Set-ExactParagraphBold {
    param($doc, [string]$exactText)

    foreach ($para in $doc.Paragraphs) {
        if ($para.Range.Text -eq ($exactText + "`r")) {
            $full = $para.Range
            $textRange = $doc.Range($full.Start, $full.End - 1)

            $flatOpc = '<?xml version="1.0" standalone="yes"?>' +
                '<?mso-application progid="Word.Document"?>' +
                '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
                '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
                '<pkg:xmlData>' +
                '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
                '<w:body><w:p><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">' + $exactText + '</w:t></w:r></w:p></w:body>' +
                '</w:document>' +
                '</pkg:xmlData></pkg:part></pkg:package>'

            $textRange.InsertXML($flatOpc)
            break
        }
    }
}

Set-ExactParagraphBold $d "todo.md file"
Set-ExactParagraphBold $d "project.yaml file"

Write-Host "Edits applied"
